$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cells ---

# Row 1 new headers
$ws.Range("G1").Value = "UserChoice"
$ws.Range("H1").Value = "Expression"

# Row 2: graphlocation simplified, add Expression
$ws.Range("E2").Value = "graph1"
$ws.Range("H2").Value = 'SUM_TIMEWIN("Adj Close", [TimeWindowStart],[TimeWindowEnd])'

# Row 3: graphlocation simplified
$ws.Range("E3").Value = "graph2"

# Row 4: graphlocation simplified
$ws.Range("E4").Value = "graph3"

# Row 5: title text tweak + graphlocation simplified
$ws.Range("A5").Value = "Verify the workbook/visualization `nfeatures for a non time `nseries data table having `nTime bucket."
$ws.Range("E5").Value = "graph4"

# --- New rows 6-9 ---

# Row 6
$ws.Range("A6").Value = "Verify the workbook/visualization `nfeatures for a non time `nseries data table having `nhaving Calculated Column for KDB connector"
$ws.Range("B6").Value = "workbook5"
$ws.Range("D6").Value = "scrollElement-visualization.VerticalBarGraph1"
$ws.Range("E6").Value = "graph5"
$ws.Range("F6").Value = "Calculated"
$ws.Range("G6").Value = "KDB"

# Row 7
$ws.Range("A7").Value = "Verify the workbook/visualization for a timeseries visualization having Calculated column added to it for MQTT connector"
$ws.Range("B7").Value = "workbook6"
$ws.Range("D7").Value = "scrollElement-visualization.LineGraph1"
$ws.Range("F7").Value = "Calculated"
$ws.Range("G7").Value = "MQTT"
$ws.Range("H7").Value = "[Adj_Close]+[Holding] "

# Row 8
$ws.Range("A8").Value = "Verify the workbook/visualization for a timeseries visualization having Ranking column added to it for KDB connector"
$ws.Range("B8").Value = "workbook7"
$ws.Range("D8").Value = "scrollElement-visualization.VerticalBarGraph1"
$ws.Range("E8").Value = "graph6"
$ws.Range("F8").Value = "Ranking"
$ws.Range("G8").Value = "KDB"

# Row 9
$ws.Range("A9").Value = "Verify the workbook/visualization features for a non time series data table having Time bucket for KDB connector"
$ws.Range("B9").Value = "workbook8"
$ws.Range("D9").Value = "scrollElement-visualization.VerticalBarGraph1"
$ws.Range("E9").Value = "graph7"
$ws.Range("F9").Value = "Time Bucket"
$ws.Range("G9").Value = "KDB"

# --- Wrap text styling on new/changed cells (matches column A convention + H2/H7) ---
$ws.Range("A6").WrapText = $true
$ws.Range("A7").WrapText = $true
$ws.Range("A8").WrapText = $true
$ws.Range("A9").WrapText = $true
$ws.Range("H2").WrapText = $true
$ws.Range("H7").WrapText = $true

# --- Column widths for new columns ---
$ws.Columns.Item(7).ColumnWidth = 12.3
$ws.Columns.Item(8).ColumnWidth = 19.3

# --- Row heights ---
$ws.Rows.Item(2).RowHeight = 86.4
$ws.Rows.Item(6).RowHeight = 72
$ws.Rows.Item(7).RowHeight = 43.2
$ws.Rows.Item(8).RowHeight = 43.2
$ws.Rows.Item(9).RowHeight = 43.2

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1

# --- Selection / scroll position ---
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("H7").Select() | Out-Null
